$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("D2").Value = "62.738.99"
$ws.Range("E2").Value = "  -8.05%  "
$ws.Range("D3").Value = "3.207.17"
$ws.Range("E3").Value = "  -10.79%  "
$ws.Range("E4").Value = "  +0.18%  "
$ws.Range("D5").Value = "'172.24"
$ws.Range("E5").Value = "  -17.07%  "
$ws.Range("D6").Value = "'505.74"
$ws.Range("E6").Value = "  -11.23%  "
$ws.Range("D7").Value = "'0.581"
$ws.Range("E7").Value = "  -5.22%  "
$ws.Range("E8").Value = "  +0.11%  "
$ws.Range("D9").Value = "3.202.77"
$ws.Range("E9").Value = "  -10.74%  "
$ws.Range("D10").Value = "'0.599"
$ws.Range("E10").Value = "  -12.44%  "
$ws.Range("D11").Value = "'55.35"
$ws.Range("E11").Value = "  -12.70%  "
$ws.Range("D12").Value = "'0.126"
$ws.Range("E12").Value = "  -14.87%  "
$ws.Range("D13").Value = "'0.0000247"
$ws.Range("E13").Value = "  -12.29%  "
$ws.Range("D14").Value = "'8.90"
$ws.Range("E14").Value = "  -14.12%  "
$ws.Range("D15").Value = "3.745.54"
$ws.Range("E15").Value = "  -10.01%  "
$ws.Range("E16").Value = "  -7.55%  "
$ws.Range("D17").Value = "3.232.36"
$ws.Range("E17").Value = "  -9.89%  "
$ws.Range("D18").Value = "62.684.80"
$ws.Range("E18").Value = "  -7.81%  "
$ws.Range("D19").Value = "'16.82"
$ws.Range("E19").Value = "  -12.74%  "
$ws.Range("D20").Value = "'10.57"
$ws.Range("E20").Value = "  -13.52%  "
$ws.Range("D21").Value = "'0.918"
$ws.Range("E21").Value = "  -13.81%  "
$ws.Range("D22").Value = "'361.40"
$ws.Range("E22").Value = "  -10.53%  "
$ws.Range("D23").Value = "'78.08"
$ws.Range("E23").Value = "  -7.94%  "
$ws.Range("B24").Value = "RenderToken"
$ws.Range("C24").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D24").Value = "'10.66"
$ws.Range("E24").Value = "  -13.92%  "
$ws.Range("B25").Value = "PancakeSwap"
$ws.Range("C25").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D25").Value = "'3.55"
$ws.Range("E25").Value = "  -14.82%  "
$ws.Range("D26").Value = "'5.91"
$ws.Range("E26").Value = "  -3.65%  "
$ws.Range("D27").Value = "'3.68"
$ws.Range("E27").Value = "  -5.06%  "
$ws.Range("D28").Value = "'2.57"
$ws.Range("E28").Value = "  -11.08%  "
$ws.Range("D29").Value = "'10.96"
$ws.Range("E29").Value = "  -12.50%  "
$ws.Range("D30").Value = "'8.04"
$ws.Range("E30").Value = "  -13.31%  "
$ws.Range("D31").Value = "'636.43"
$ws.Range("E31").Value = "  -8.15%  "
$ws.Range("D32").Value = "'27.58"
$ws.Range("E32").Value = "  -12.70%  "
$ws.Range("D33").Value = "'6.45"
$ws.Range("E33").Value = "  -15.21%  "
$ws.Range("D34").Value = "'10.84"
$ws.Range("E34").Value = "  -10.90%  "
$ws.Range("D35").Value = "'58.30"
$ws.Range("E35").Value = "  -8.23%  "
$ws.Range("B36").Value = "Dai"
$ws.Range("C36").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D36").Value = "'0.999"
$ws.Range("E36").Value = "  -0.14%  "
$ws.Range("B37").Value = "Hedera"
$ws.Range("C37").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D37").Value = "'0.100"
$ws.Range("E37").Value = "  -11.22%  "
$ws.Range("D38").Value = "'34.86"
$ws.Range("E38").Value = "  -16.29%  "
$ws.Range("D39").Value = "'0.367"
$ws.Range("E39").Value = "  -10.60%  "
$ws.Range("D40").Value = "'1.00"
$ws.Range("E40").Value = "  +0.27%  "
$ws.Range("B41").Value = "Kaspa"
$ws.Range("C41").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D41").Value = "'0.121"
$ws.Range("E41").Value = "  -9.01%  "
$ws.Range("B42").Value = "Maker"
$ws.Range("C42").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D42").Value = "2.826.52"
$ws.Range("E42").Value = "  -10.57%  "
$ws.Range("D43").Value = "0.0₃0631"
$ws.Range("E43").Value = "  -16.88%  "
$ws.Range("D44").Value = "'2.59"
$ws.Range("E44").Value = "  -20.44%  "
$ws.Range("D45").Value = "'2.55"
$ws.Range("E45").Value = "  -7.62%  "
$ws.Range("D46").Value = "'2.28"
$ws.Range("E46").Value = "  -15.07%  "
$ws.Range("D47").Value = "'2.73"
$ws.Range("E47").Value = "  +0.83%  "
$ws.Range("D48").Value = "'0.0371"
$ws.Range("E48").Value = "  -10.30%  "
$ws.Range("B49").Value = "ApeXProtocol"
$ws.Range("C49").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D49").Value = "'2.88"
$ws.Range("E49").Value = "  -8.39%  "
$ws.Range("B50").Value = "Stellar"
$ws.Range("C50").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D50").Value = "'0.120"
$ws.Range("E50").Value = "  -8.04%  "
$ws.Range("D51").Value = "'129.34"
$ws.Range("E51").Value = "  -6.74%  "
